$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New applicant row (17) mirrors the layout/formatting of the row above it (16).
$ws.Range("A16:M16").Copy()
$ws.Range("A17").PasteSpecial(-4122)

# Populate the new record's data (order matches the shared-string table
# growth of the source workbook: nama, email, ttl, alamat are newly
# introduced strings, so they are written in that order).
$ws.Range("D17").Value = "Ramlan"
$ws.Range("A17").Value = "ramlankhenzo79@gmail.com"
$ws.Range("E17").Value = "ramlankhenzo79@gmail.com"
$ws.Range("F17").Value = "Bone 06 0ktober 2006"
$ws.Range("G17").Value = "Sulawesi tangah kabupaten poso"
$ws.Range("B17").Value = 2025
$ws.Range("C17").Value = "profil_siswa.html"
$ws.Range("H17").Value = 85810936860
$ws.Range("I17").Value = 2025
$ws.Range("J17").Value = "Putih"
$ws.Range("K17").Value = "Panding"
$ws.Range("L17").Value = "Belum tersedia (Status Pending)"
$ws.Range("M17").Value = "default"

# Restore the previously active selection to D10.
$ws.Range("D10").Select()
